# 'finished' skeleton for clueGame package
# - Remove the two unused placeholder sheets (Sheet2, Sheet3)
# - Rename the remaining sheet (Sheet1) to "ClueBoard"
# - Move the active selection from Y29 to Y11

$wb = $excel.ActiveWorkbook

# Avoid the "you are about to delete ..." confirmation prompt when
# removing worksheets.
$excel.DisplayAlerts = $false

$wb.Worksheets.Item("Sheet3").Delete()
$wb.Worksheets.Item("Sheet2").Delete()

$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "ClueBoard"

$ws.Range("Y11").Select()
